$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 4 - reorder "Recorded By" list
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 5 - reorder "Recorded By" list
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 6 - reorder "Recorded By" list, Recorded Sessions count +1
$ws.Range("G6").Value = "Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("L6").Value = 21

# Row 7 - reorder "Recorded By" list, Missing Sessions count -1
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg"
$ws.Range("L7").Value = 2

# Row 9 - reorder "Recorded By" list, Coverage % recalculated
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "72.4%"

# Row 10 - Average Attendance % recalculated
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "25.4%"

# Row 12 - reorder "Recorded By" list
$ws.Range("G12").Value = "Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"

# Row 15 - Group Statistics summary row, mirrors the recorded/missing/coverage updates
$ws.Range("O15").Value = 21
$ws.Range("P15").Value = 2
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "72.4%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "25.4%"

# Row 17 - reorder "Recorded By" list
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# Row 21 - session now recorded (was "Not Recorded"/pink style, becomes normal "Recorded" style)
$ws.Range("A2:I2").Copy()
$ws.Range("A21:I21").PasteSpecial(-4122)
$ws.Range("G21").Value = "esraa.sami@med.asu.edu.eg"
$ws.Range("H21").Value = "6/251"
$ws.Range("I21").Value = "Recorded"

# Row 28 - reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# Row 30 - reorder "Recorded By" list
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
